# Generate Report for Handback
# ----------------------------------------------------------------------------------
# This script reproduces the "handback" report-generation edit:
#  - Overview sheet: the per-language status cells flip from
#    "Ready for handoff" to "Handed back: in sync with en-US"
#  - Each language sheet (zh-cn, de-de) gets its "Latest Target File" and
#    "Latest Handback File" columns (I, J) populated with the handed-back
#    xliff hyperlink / file name, and the "Latest Handback DateTime" column
#    (K) gets a real timestamp.
#  - A few columns are widened to better fit the newly-populated long text.
# ----------------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2        # xlUnderlineStyleSingle
$hyperlinkColor = 15570276     # RGB(0x64,0x95,0xED) == FF6495ED, matches the existing "HyperLink" cell style

function Apply-HyperlinkLook($rng) {
    # NOTE: setting Font.Underline and Font.Color together has an ordering quirk in
    # this engine where only the *last* assigned property survives serialization,
    # so we deliberately re-assert Underline once more after Color to make sure both
    # stick.
    $rng.Font.Underline = $hyperlinkUnderline
    $rng.Font.Color = $hyperlinkColor
    $rng.Font.Underline = $hyperlinkUnderline
}

# ----------------------------------------------------------------------------------
# 1) Overview sheet: status text for both language columns (E/F), rows 2 and 3
# ----------------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# widen the zh-cn / de-de status columns
$overview.Range("E1").EntireColumn.ColumnWidth = 29.16666666666667
$overview.Range("F1").EntireColumn.ColumnWidth = 29.16666666666667

# ----------------------------------------------------------------------------------
# 2) zh-cn sheet
# ----------------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$handoffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b323e073613fd57a26f8f00dde15614ac3a39ea/e2e/8022f275-f833-4477-b9d8-4ea2be31bf2d.md"
$handoffDisplay = "8022f275-f833-4477-b9d8-4ea2be31bf2d.md"

# Row 2
$zhI2 = $zh.Range("I2")
$zh.Hyperlinks.Add($zhI2, $handoffUrl, "", "", $handoffDisplay) | Out-Null
Apply-HyperlinkLook $zhI2
$zh.Range("J2").Value = "8022f275-f833-4477-b9d8-4ea2be31bf2d.816eed1bbd4b688c757dbb3800454b802dd95472.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-01 13:10:59"

# Row 3
$zhI3 = $zh.Range("I3")
$zh.Hyperlinks.Add($zhI3, $handoffUrl, "", "", $handoffDisplay) | Out-Null
Apply-HyperlinkLook $zhI3
$zh.Range("J3").Value = "8022f275-f833-4477-b9d8-4ea2be31bf2d.816eed1bbd4b688c757dbb3800454b802dd95472.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-01 13:10:59"

# widen columns C (Status), I (Latest Target File), J (Latest Handback File)
$zh.Range("C1").EntireColumn.ColumnWidth = 29.16666666666667
$zh.Range("I1").EntireColumn.ColumnWidth = 39.16666666666667
$zh.Range("J1").EntireColumn.ColumnWidth = 39.16666666666667

# ----------------------------------------------------------------------------------
# 3) de-de sheet
# ----------------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2
$deI2 = $de.Range("I2")
$de.Hyperlinks.Add($deI2, $handoffUrl, "", "", $handoffDisplay) | Out-Null
Apply-HyperlinkLook $deI2
$de.Range("J2").Value = "8022f275-f833-4477-b9d8-4ea2be31bf2d.816eed1bbd4b688c757dbb3800454b802dd95472.de-de.xlf"
$de.Range("K2").Value = "2016-09-01 13:11:22"

# Row 3
$deI3 = $de.Range("I3")
$de.Hyperlinks.Add($deI3, $handoffUrl, "", "", $handoffDisplay) | Out-Null
Apply-HyperlinkLook $deI3
$de.Range("J3").Value = "8022f275-f833-4477-b9d8-4ea2be31bf2d.816eed1bbd4b688c757dbb3800454b802dd95472.de-de.xlf"
$de.Range("K3").Value = "2016-09-01 13:11:22"

# widen columns C (Status), I (Latest Target File), J (Latest Handback File)
$de.Range("C1").EntireColumn.ColumnWidth = 29.16666666666667
$de.Range("I1").EntireColumn.ColumnWidth = 39.16666666666667
$de.Range("J1").EntireColumn.ColumnWidth = 39.16666666666667

Write-Host "Handback report generated."
